$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("D16").Value = 44186
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 21000
$ws.Range("L16").Value = 22000
$ws.Range("M16").Value = 21500
$ws.Range("N16").Value = "`$/saco 25 kilos"
$ws.Range("O16").Value = "Provincia de Diguillín"
$ws.Range("P16").Value = 860
$ws.Range("Q16").Value = 25

# Row 17
$ws.Range("D17").Value = 44215
$ws.Range("J17").Value = 42
$ws.Range("K17").Value = 19000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 19524
$ws.Range("O17").Value = "Región de La Araucanía"
$ws.Range("P17").Value = 781

# Row 18
$ws.Range("D18").Value = 44544
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 13000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13500
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 540

# Row 19
$ws.Range("D19").Value = 44530
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14500
$ws.Range("P19").Value = 580

# Row 20
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 13000
$ws.Range("M20").Value = 12500
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 500

# Row 21
$ws.Range("D21").Value = 44557
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("O21").Value = "Provincia de Diguillín"
$ws.Range("P21").Value = 700

# Row 22
$ws.Range("D22").Value = 44537
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 13000
$ws.Range("L22").Value = 14000
$ws.Range("M22").Value = 13500
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 540

# Row 23
$ws.Range("D23").Value = 44546
$ws.Range("J23").Value = 60
$ws.Range("O23").Value = "Provincia de Diguillín"

# Row 24
$ws.Range("D24").Value = 44162
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 18500
$ws.Range("L24").Value = 19000
$ws.Range("M24").Value = 18820
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 753

# Row 25
$ws.Range("D25").Value = 44554
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 15500
$ws.Range("O25").Value = "Provincia de Diguillín"
$ws.Range("P25").Value = 620

# Row 26
$ws.Range("D26").Value = 44166
$ws.Range("J26").Value = 48
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17479
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 699

# Row 27
$ws.Range("D27").Value = 44174
$ws.Range("H27").Value = "Perfection"
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = 19000
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = 19500
$ws.Range("O27").Value = "Región de Ñuble"
$ws.Range("P27").Value = 780

# Row 28
$ws.Range("D28").Value = 44273
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = 20000
$ws.Range("L28").Value = 22000
$ws.Range("M28").Value = 21091
$ws.Range("O28").Value = "Región de La Araucanía"
$ws.Range("P28").Value = 844

# Row 29
$ws.Range("D29").Value = 44525
$ws.Range("J29").Value = 80
$ws.Range("K29").Value = 12000
$ws.Range("L29").Value = 13000
$ws.Range("M29").Value = 12500
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 500

# Row 30
$ws.Range("D30").Value = 44165
$ws.Range("J30").Value = 42
$ws.Range("K30").Value = 18000
$ws.Range("L30").Value = 19000
$ws.Range("M30").Value = 18595
$ws.Range("P30").Value = 744

# Row 31
$ws.Range("D31").Value = 44516
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 16000
$ws.Range("M31").Value = 15500
$ws.Range("P31").Value = 620

# Row 32
$ws.Range("D32").Value = 44567
$ws.Range("J32").Value = 60
$ws.Range("K32").Value = 18000
$ws.Range("L32").Value = 19000
$ws.Range("M32").Value = 18500
$ws.Range("O32").Value = "Provincia de Diguillín"
$ws.Range("P32").Value = 740

# Row 48
$ws.Range("D48").Value = 44568
$ws.Range("J48").Value = 120
$ws.Range("K48").Value = 24000
$ws.Range("L48").Value = 25000
$ws.Range("M48").Value = 24500
$ws.Range("O48").Value = "Provincia de Diguillín"
$ws.Range("P48").Value = 980

# Row 49
$ws.Range("D49").Value = 44518
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 16000
$ws.Range("M49").Value = 15500
$ws.Range("P49").Value = 620

# Row 50
$ws.Range("A50").Value = 7
$ws.Range("B50").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value = "Ñuble"
$ws.Range("D50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D50").Value = 44540
$ws.Range("E50").Value = 16
$ws.Range("F50").Value = 100112022
$ws.Range("G50").Value = "Arveja Verde"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 13000
$ws.Range("L50").Value = 14000
$ws.Range("M50").Value = 13500
$ws.Range("N50").Value = "`$/saco 25 kilos"
$ws.Range("O50").Value = "Región del Maule"
$ws.Range("P50").Value = 540
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
